$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("final_fail")
$ws1.Range("A9").Value = 'Start of Session 1 (%)'
$ws1.Range("A10").Value = 'Clicks per day'
$ws1.Range("B10").Value = $false
$ws1.Range("D10").Value = $true
$ws1.Range("A11").Value = 'Number of days'
$ws1.Range("A12").Value = 'Resources viewed'
$ws1.Range("D12").Value = $true
$ws1.Range("E12").Value = $false
$ws1.Range("A13").Value = 'Number of clicks'
$ws1.Range("D13").Value = $false
$ws1.Range("J13").Value = 4
$ws1.Range("A14").Value = 'Clicks per session'
$ws1.Range("A15").Value = 'Clicks on folder'
$ws1.Range("B15").Value = $true
$ws1.Range("F15").Value = $false
$ws1.Range("A16").Value = 'Total time online (min)'
$ws1.Range("D16").Value = $false
$ws1.Range("F16").Value = $true
$ws1.Range("A18").Value = 'Start of Session 2 (%)'
$ws1.Range("A19").Value = 'Start of Session 3 (%)'
$ws1.Range("A22").Value = 'Start of Session 7 (%)'
$ws1.Range("A23").Value = 'Assignments viewed'
$ws1.Range("B23").Value = $true
$ws1.Range("F23").Value = $false
$ws1.Range("A24").Value = 'Start of Session 4 (%)'
$ws1.Range("B24").Value = $false
$ws1.Range("F24").Value = $true
$ws1.Range("A28").Value = 'Forum posts'
$ws1.Range("C28").Value = $false
$ws1.Range("D28").Value = $true
$ws1.Range("A29").Value = 'Start of Session 10 (%)'
$ws1.Range("C29").Value = $true
$ws1.Range("D29").Value = $false
$ws1.Range("A30").Value = 'Quizzes started'
$ws1.Range("C30").Value = $true
$ws1.Range("D30").Value = $false
$ws1.Range("A31").Value = 'Discussions viewed'
$ws1.Range("C31").Value = $false
$ws1.Range("D31").Value = $true
$ws1.Range("A32").Value = 'Assignments submitted'
$ws1.Range("A33").Value = 'Clicks on course'
$ws1.Range("A34").Value = 'Number of sessions'
$ws1.Range("A35").Value = 'Clicks on forum'
$ws1.Range("C35").Value = $true
$ws1.Range("J35").Value = 2
$ws1.Range("A36").Value = 'Start of Session 8 (%)'
$ws1.Range("A37").Value = 'Start of Session 9 (%)'

$ws2 = $wb.Worksheets.Item("final_gifted")
$ws2.Range("A2").Value = 'Average grade of assignments'
$ws2.Range("A3").Value = 'Largest period of inactivity (h)'
$ws2.Range("A4").Value = 'Average session duration (min)'
$ws2.Range("A5").Value = 'Start of Session 1 (%)'
$ws2.Range("A6").Value = 'Clicks (% of course total)'
$ws2.Range("A7").Value = 'Resources viewed'
$ws2.Range("D7").Value = $false
$ws2.Range("E7").Value = $true
$ws2.Range("A8").Value = 'Clicks per session'
$ws2.Range("D8").Value = $true
$ws2.Range("E8").Value = $false
$ws2.Range("A9").Value = 'Total time online (min)'
$ws2.Range("C9").Value = $false
$ws2.Range("D9").Value = $true
$ws2.Range("A10").Value = 'On/off campus click ratio'
$ws2.Range("D10").Value = $false
$ws2.Range("E10").Value = $true
$ws2.Range("A11").Value = 'Days with no interaction'
$ws2.Range("E11").Value = $true
$ws2.Range("J11").Value = 5
$ws2.Range("A12").Value = 'Assignments viewed'
$ws2.Range("C12").Value = $true
$ws2.Range("E12").Value = $false
$ws2.Range("A13").Value = 'Assignments submitted'
$ws2.Range("C13").Value = $false
$ws2.Range("D13").Value = $true
$ws2.Range("F13").Value = $false
$ws2.Range("J13").Value = 3
$ws2.Range("A14").Value = 'Number of days'
$ws2.Range("A15").Value = 'Start of Session 6 (%)'
$ws2.Range("D15").Value = $false
$ws2.Range("F15").Value = $true
$ws2.Range("A16").Value = 'Start of Session 4 (%)'
$ws2.Range("A19").Value = 'Clicks per day'
$ws2.Range("A20").Value = 'Submissions (% of course total)'
$ws2.Range("D20").Value = $false
$ws2.Range("F20").Value = $true
$ws2.Range("A21").Value = 'Files downloaded'
$ws2.Range("A22").Value = 'Quizzes started'
$ws2.Range("D22").Value = $true
$ws2.Range("F22").Value = $false
$ws2.Range("A23").Value = 'Number of clicks'
$ws2.Range("D23").Value = $false
$ws2.Range("F23").Value = $true
$ws2.Range("A24").Value = 'Clicks on folder'
$ws2.Range("B24").Value = $false
$ws2.Range("D24").Value = $true
$ws2.Range("A25").Value = 'Start of Session 8 (%)'
$ws2.Range("B25").Value = $true
$ws2.Range("F25").Value = $false
$ws2.Range("A26").Value = 'Start of Session 7 (%)'
$ws2.Range("A27").Value = 'Start of Session 5 (%)'
$ws2.Range("A32").Value = 'Clicks on course'
$ws2.Range("A34").Value = 'Number of sessions'
$ws2.Range("A35").Value = 'Clicks on forum'
$ws2.Range("A36").Value = 'Start of Session 10 (%)'
$ws2.Range("A37").Value = 'Start of Session 9 (%)'
